$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-6 (regen save_data to use K instead of Strike#)
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 2
